# Insert two new data rows immediately above the current row 724
# (i.e. right before the "2021-07-20 / Fukumoto" pair), pushing the
# existing rows 724:791 down to 726:793. The two newly inserted rows
# are populated with a new "Lane Late" price observation (Primera and
# Segunda quality), dated 2023-08-28 (serial 45166).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 724; formatting (incl. the date
# number-format on column D) is inherited from the row above, same as
# interactive Excel row-insert.
$ws.Rows.Item(724).Insert()
$ws.Rows.Item(724).Insert()

# --- New row 724: Lane Late / Primera -------------------------------
$ws.Cells.Item(724, 1).Value2  = 7
$ws.Cells.Item(724, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(724, 3).Value2  = "Ñuble"
$ws.Cells.Item(724, 4).Value2  = 45166
$ws.Cells.Item(724, 5).Value2  = 16
$ws.Cells.Item(724, 6).Value2  = "Fruta"
$ws.Cells.Item(724, 7).Value2  = 100102
$ws.Cells.Item(724, 8).Value2  = "Cítricos"
$ws.Cells.Item(724, 9).Value2  = 100102005
$ws.Cells.Item(724, 10).Value2 = "Naranja"
$ws.Cells.Item(724, 11).Value2 = "Lane Late"
$ws.Cells.Item(724, 12).Value2 = "Primera"
$ws.Cells.Item(724, 13).Value2 = 120
$ws.Cells.Item(724, 14).Value2 = 10000
$ws.Cells.Item(724, 15).Value2 = 10000
$ws.Cells.Item(724, 16).Value2 = 10000
$ws.Cells.Item(724, 17).Value2 = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(724, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(724, 19).Value2 = 667
$ws.Cells.Item(724, 20).Value2 = 15

# --- New row 725: Lane Late / Segunda -------------------------------
$ws.Cells.Item(725, 1).Value2  = 7
$ws.Cells.Item(725, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(725, 3).Value2  = "Ñuble"
$ws.Cells.Item(725, 4).Value2  = 45166
$ws.Cells.Item(725, 5).Value2  = 16
$ws.Cells.Item(725, 6).Value2  = "Fruta"
$ws.Cells.Item(725, 7).Value2  = 100102
$ws.Cells.Item(725, 8).Value2  = "Cítricos"
$ws.Cells.Item(725, 9).Value2  = 100102005
$ws.Cells.Item(725, 10).Value2 = "Naranja"
$ws.Cells.Item(725, 11).Value2 = "Lane Late"
$ws.Cells.Item(725, 12).Value2 = "Segunda"
$ws.Cells.Item(725, 13).Value2 = 120
$ws.Cells.Item(725, 14).Value2 = 8000
$ws.Cells.Item(725, 15).Value2 = 8000
$ws.Cells.Item(725, 16).Value2 = 8000
$ws.Cells.Item(725, 17).Value2 = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(725, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(725, 19).Value2 = 533
$ws.Cells.Item(725, 20).Value2 = 15
